# Fruta / hortaliza, semanal
# Insert a new weekly record at row 35, pushing the existing rows 35-55
# down to 36-56 (dimension grows from A1:T55 to A1:T56).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 35; this shifts rows 35..55 -> 36..56
$ws.Rows(35).Insert()

# Populate the newly inserted row 35 with the new weekly record.
$ws.Range("A35").Value = 10
$ws.Range("B35").Value = "Vega Modelo de Temuco"
$ws.Range("C35").Value = "La Araucanía"
$ws.Range("D35").Value = 44567
$ws.Range("D35").NumberFormat = $ws.Range("D36").NumberFormat
$ws.Range("E35").Value = 9
$ws.Range("F35").Value = "Fruta"
$ws.Range("G35").Value = 100103
$ws.Range("H35").Value = "Frutos de hueso (carozo)"
$ws.Range("I35").Value = 100103003
$ws.Range("J35").Value = "Damasco"
$ws.Range("K35").Value = "Modesto"
$ws.Range("L35").Value = "Primera"
$ws.Range("M35").Value = 25
$ws.Range("N35").Value = 15000
$ws.Range("O35").Value = 15000
$ws.Range("P35").Value = 15000
$ws.Range("Q35").Value = "$/bandeja 10 kilos"
$ws.Range("R35").Value = "Provincia de Quillota"
$ws.Range("S35").Value = 1500
$ws.Range("T35").Value = 10
